$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("pie")

# Replace the numeric genotype counts with the actual allele letters
$ws.Range("B2").Value = "A"
$ws.Range("C2").Value = "C"

$ws.Range("B3").Value = "A"
$ws.Range("C3").Value = "C"

$ws.Range("B4").Value = "A"
$ws.Range("C4").Value = "T"

$ws.Range("B5").Value = "C"
$ws.Range("C5").Value = "T"

# Add a new row of explanation data
$ws.Range("B6").Value = "T"
$ws.Range("C6").Value = "T"

# Update the selection to match the author's final cursor position
[void]$ws.Range("B12").Select()
